$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as TEXT (matching the source inlineStr cells),
# without leaving a residual NumberFormat style on the cell.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue 'D2' '26.300.35'
Set-TextValue 'E2' '  -3.29%  '
Set-TextValue 'D3' '1.791.85'
Set-TextValue 'E3' '  -3.06%  '
Set-TextValue 'D4' '1.008'
Set-TextValue 'E4' '  +0.70%  '
Set-TextValue 'D5' '1.006'
Set-TextValue 'E5' '  +0.50%  '
Set-TextValue 'D6' '306.73'
Set-TextValue 'E6' '  -2.14%  '
Set-TextValue 'D7' '0.4546'
Set-TextValue 'E7' '  -1.76%  '
Set-TextValue 'D8' '0.3618'
Set-TextValue 'E8' '  -2.09%  '
Set-TextValue 'D9' '0.07064'
Set-TextValue 'E9' '  -2.81%  '
Set-TextValue 'D10' '0.8706'
Set-TextValue 'E10' '  -1.83%  '
Set-TextValue 'D11' '0.07778'
Set-TextValue 'E11' '  -0.45%  '
Set-TextValue 'D12' '19.33'
Set-TextValue 'E12' '  -2.99%  '
Set-TextValue 'D13' '1.746.84'
Set-TextValue 'E13' '  -5.92%  '
Set-TextValue 'D14' '5.259'
Set-TextValue 'E14' '  -2.42%  '
Set-TextValue 'D15' '6.313'
Set-TextValue 'D16' '84.51'
Set-TextValue 'E16' '  -7.72%  '
Set-TextValue 'E17' '  +0.69%  '
Set-TextValue 'D18' '0.000008492'
Set-TextValue 'E18' '  -4.05%  '
Set-TextValue 'E19' '  +0.60%  '
Set-TextValue 'D20' '26.358.69'
Set-TextValue 'E20' '  -3.15%  '
Set-TextValue 'D21' '14.13'
Set-TextValue 'E21' '  -3.59%  '
Set-TextValue 'D22' '4.973'
Set-TextValue 'E22' '  -1.60%  '
Set-TextValue 'E23' '  -0.56%  '
Set-TextValue 'D24' '1.975.90'
Set-TextValue 'E24' '  -6.75%  '
Set-TextValue 'D25' '1.974'
Set-TextValue 'E25' '  -3.13%  '
Set-TextValue 'D26' '151.92'
Set-TextValue 'E26' '  +0.38%  '
Set-TextValue 'D27' '17.78'
Set-TextValue 'E27' '  -3.10%  '
Set-TextValue 'D28' '2.043'
Set-TextValue 'E28' '  +0.51%  '
Set-TextValue 'D29' '112.27'
Set-TextValue 'E29' '  -2.79%  '
Set-TextValue 'D30' '4.827'
Set-TextValue 'E30' '  -3.80%  '
Set-TextValue 'D31' '0.08647'
Set-TextValue 'E31' '  -2.08%  '
Set-TextValue 'D32' '3.025'
Set-TextValue 'E32' '  -3.95%  '
Set-TextValue 'D33' '4.436'
Set-TextValue 'E33' '  -1.56%  '
Set-TextValue 'B34' 'RenderToken'
Set-TextValue 'C34' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D34' '2.682'
Set-TextValue 'E34' '  -1.33%  '
Set-TextValue 'B35' 'ImmutableX'
Set-TextValue 'C35' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D35' '0.7141'
Set-TextValue 'E35' '  -9.43%  '
Set-TextValue 'D36' '1.106'
Set-TextValue 'E36' '  -3.90%  '
Set-TextValue 'D37' '1.005'
Set-TextValue 'E37' '  +0.53%  '
Set-TextValue 'D38' '1.076'
Set-TextValue 'E38' '  -2.45%  '
Set-TextValue 'D39' '0.01939'
Set-TextValue 'E39' '  -0.29%  '
Set-TextValue 'D40' '0.05080'
Set-TextValue 'E40' '  -2.75%  '
Set-TextValue 'D41' '2.858'
Set-TextValue 'E41' '  -3.14%  '
Set-TextValue 'D42' '6.879'
Set-TextValue 'E42' '  -2.31%  '
Set-TextValue 'D43' '0.4915'
Set-TextValue 'E43' '  -2.49%  '
Set-TextValue 'D44' '0.1513'
Set-TextValue 'E44' '  -6.17%  '
Set-TextValue 'D45' '7.963'
Set-TextValue 'E45' '  -6.28%  '
Set-TextValue 'D46' '1.006'
Set-TextValue 'E46' '  +0.54%  '
Set-TextValue 'D47' '0.4568'
Set-TextValue 'E47' '  -3.91%  '
Set-TextValue 'D48' '9.855'
Set-TextValue 'E48' '  -4.64%  '
Set-TextValue 'D49' '99.68'
Set-TextValue 'E49' '  -3.11%  '
Set-TextValue 'D50' '1.579'
Set-TextValue 'E50' '  -3.54%  '
Set-TextValue 'D51' '0.05948'
Set-TextValue 'E51' '  -3.96%  '
